# Auto-generated edit script: updates H:N market-profit columns per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1493.4
$ws.Range("J19").Value = 2000
$ws.Range("L19").Value = 2000
$ws.Range("N19").Value = -2350
$ws.Range("H53").Value = 476.33334
$ws.Range("I53").Value = 398.8
$ws.Range("K53").Value = 398.8
$ws.Range("M53").Value = 238.2
$ws.Range("H70").Value = 1833.3334
$ws.Range("I70").Value = 3000
$ws.Range("K70").Value = 9000
$ws.Range("M70").Value = -8730
$ws.Range("H73").Value = 1833.3334
$ws.Range("I73").Value = 3000
$ws.Range("K73").Value = 9000
$ws.Range("M73").Value = -8064
$ws.Range("H111").Value = 369.75
$ws.Range("I111").Value = 224.5
$ws.Range("K111").Value = 673.5
$ws.Range("M111").Value = 2393.5
$ws.Range("H112").Value = 2696.5625
$ws.Range("J112").Value = 2696.5625
$ws.Range("L112").Value = 8089.6875
$ws.Range("N112").Value = -10305.6875
$ws.Range("H137").Value = 2269.75
$ws.Range("I137").Value = 2269.75
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 6809.25
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -4259.25
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("H45").Value = 32082.5
$ws.Range("I45").Value = 3123.75
$ws.Range("J45").Value = 90000
$ws.Range("K45").Value = 3123.75
$ws.Range("L45").Value = 90000
$ws.Range("M45").Value = -2746.75
$ws.Range("N45").Value = -90754
$ws.Range("H63").Value = 11580.777
$ws.Range("J63").Value = 12103.667
$ws.Range("L63").Value = 12103.667
$ws.Range("N63").Value = -13475.667
$ws.Range("H66").Value = 11580.777
$ws.Range("J66").Value = 12103.667
$ws.Range("L66").Value = 60518.335
$ws.Range("N66").Value = -67382.33499999999
$ws.Range("H74").Value = 4061.125
$ws.Range("I74").Value = 4251.9
$ws.Range("J74").Value = 1199.5
$ws.Range("K74").Value = 4251.9
$ws.Range("L74").Value = 1199.5
$ws.Range("M74").Value = -3377.9
$ws.Range("N74").Value = -2947.5
$ws.Range("H77").Value = 4061.125
$ws.Range("I77").Value = 4251.9
$ws.Range("J77").Value = 1199.5
$ws.Range("K77").Value = 21259.5
$ws.Range("L77").Value = 5997.5
$ws.Range("M77").Value = -16891.5
$ws.Range("N77").Value = -14733.5
$ws.Range("H80").Value = 33621.2
$ws.Range("I80").Value = 21499
$ws.Range("J80").Value = 82110
$ws.Range("K80").Value = 21499
$ws.Range("L80").Value = 82110
$ws.Range("M80").Value = -20501
$ws.Range("N80").Value = -84106
$ws.Range("H83").Value = 33621.2
$ws.Range("I83").Value = 21499
$ws.Range("J83").Value = 82110
$ws.Range("K83").Value = 64497
$ws.Range("L83").Value = 246330
$ws.Range("M83").Value = -59505
$ws.Range("N83").Value = -256314
$ws.Range("H110").Value = 3877.182
$ws.Range("I110").Value = 1805.5555
$ws.Range("K110").Value = 1805.5555
$ws.Range("M110").Value = 239.4445000000001
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 4029.2
$ws.Range("I132").Value = 2970.5715
$ws.Range("K132").Value = 8911.7145
$ws.Range("M132").Value = -6381.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("H80").Value = 416.4
$ws.Range("I80").Value = 124
$ws.Range("J80").Value = 541.7143
$ws.Range("K80").Value = 124
$ws.Range("L80").Value = 541.7143
$ws.Range("M80").Value = 874
$ws.Range("N80").Value = -2537.7143
$ws.Range("H83").Value = 416.4
$ws.Range("I83").Value = 124
$ws.Range("J83").Value = 541.7143
$ws.Range("K83").Value = 620
$ws.Range("L83").Value = 2708.5715
$ws.Range("M83").Value = 4372
$ws.Range("N83").Value = -12692.5715
$ws.Range("H134").Value = 5912.143
$ws.Range("I134").Value = 2730.8333
$ws.Range("J134").Value = 25000
$ws.Range("K134").Value = 8192.499899999999
$ws.Range("L134").Value = 75000
$ws.Range("M134").Value = -5657.499899999999
$ws.Range("N134").Value = -80070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H58").Value = 3353.913
$ws.Range("I58").Value = 3344.7368
$ws.Range("J58").Value = 3397.5
$ws.Range("K58").Value = 3344.7368
$ws.Range("L58").Value = 3397.5
$ws.Range("M58").Value = -3141.7368
$ws.Range("N58").Value = -3803.5
$ws.Range("H68").Value = 78645
$ws.Range("J68").Value = 78645
$ws.Range("L68").Value = 78645
$ws.Range("N68").Value = -80143
$ws.Range("H71").Value = 78645
$ws.Range("J71").Value = 78645
$ws.Range("L71").Value = 235935
$ws.Range("N71").Value = -243423
$ws.Range("H86").Value = 23327.5
$ws.Range("I86").Value = 23327.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 23327.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -22204.5
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 23327.5
$ws.Range("I89").Value = 23327.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 116637.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -111021.5
$ws.Range("N89").ClearContents()
$ws.Range("H99").Value = 2653
$ws.Range("I99").Value = 2653
$ws.Range("K99").Value = 2653
$ws.Range("M99").Value = -1155
$ws.Range("H107").Value = 2470.3333
$ws.Range("I107").Value = 4000
$ws.Range("J107").Value = 1705.5
$ws.Range("K107").Value = 4000
$ws.Range("L107").Value = 1705.5
$ws.Range("M107").Value = -2080
$ws.Range("N107").Value = -5545.5
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 2653
$ws.Range("I126").Value = 2653
$ws.Range("K126").Value = 7959
$ws.Range("M126").Value = -5489
$ws.Range("H136").Value = 3353.913
$ws.Range("I136").Value = 3344.7368
$ws.Range("J136").Value = 3397.5
$ws.Range("K136").Value = 10034.2104
$ws.Range("L136").Value = 10192.5
$ws.Range("M136").Value = -7484.2104
$ws.Range("N136").Value = -15292.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 56.272728
$ws.Range("I2").Value = 50.8
$ws.Range("K2").Value = 304.8
$ws.Range("M2").Value = -191.8
$ws.Range("H23").Value = 146
$ws.Range("I23").Value = 24.5
$ws.Range("J23").Value = 267.5
$ws.Range("K23").Value = 73.5
$ws.Range("L23").Value = 802.5
$ws.Range("M23").Value = 161.5
$ws.Range("N23").Value = -1272.5
$ws.Range("H29").Value = 1046.2858
$ws.Range("J29").Value = 1174
$ws.Range("L29").Value = 3522
$ws.Range("N29").Value = -4076
$ws.Range("H34").Value = 4301.0625
$ws.Range("J34").Value = 5250.25
$ws.Range("L34").Value = 15750.75
$ws.Range("N34").Value = -15918.75
$ws.Range("H38").Value = 209.75
$ws.Range("J38").Value = 199.66667
$ws.Range("L38").Value = 599.00001
$ws.Range("N38").Value = -1293.00001
$ws.Range("H92").Value = 499.33334
$ws.Range("I92").Value = 399.2
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 1197.6
$ws.Range("L92").Value = 3000
$ws.Range("M92").Value = 50.40000000000009
$ws.Range("N92").Value = -5496
$ws.Range("H107").Value = 337.5
$ws.Range("I107").Value = 100
$ws.Range("K107").Value = 300
$ws.Range("M107").Value = 1620
$ws.Range("H131").Value = 2087.3428
$ws.Range("I131").Value = 2472.5
$ws.Range("K131").Value = 7417.5
$ws.Range("M131").Value = -2377.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 162.57143
$ws.Range("I2").Value = 173
$ws.Range("K2").Value = 173
$ws.Range("M2").Value = -60
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H132").Value = 4499.5
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11149.714
$ws.Range("I22").Value = 15782.667
$ws.Range("K22").Value = 15782.667
$ws.Range("M22").Value = -15487.667
$ws.Range("H27").Value = 11149.714
$ws.Range("I27").Value = 15782.667
$ws.Range("K27").Value = 15782.667
$ws.Range("M27").Value = -15675.667
$ws.Range("H55").Value = 2690
$ws.Range("I55").Value = 3900
$ws.Range("K55").Value = 3900
$ws.Range("M55").Value = -3727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 41662.668
$ws.Range("J74").Value = 39994.5
$ws.Range("L74").Value = 39994.5
$ws.Range("N74").Value = -41866.5
$ws.Range("H77").Value = 41662.668
$ws.Range("J77").Value = 39994.5
$ws.Range("L77").Value = 119983.5
$ws.Range("N77").Value = -129343.5
